$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain stored as text, matching the source data feed format

$ws.Range("D2").Value = "29.865.47"
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").Value = "1.893.90"
$ws.Range("E3").Value = "  +0.32%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7816"
$ws.Range("E5").Value = "  +0.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.87"
$ws.Range("E6").Value = "  +1.00%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3132"
$ws.Range("E8").Value = "  -0.67%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.74"
$ws.Range("E9").Value = "  +1.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07338"
$ws.Range("E10").Value = "  +5.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08095"
$ws.Range("E11").Value = "  +0.81%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7718"
$ws.Range("E12").Value = "  +1.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.495"
$ws.Range("E13").Value = "  +4.62%  "

$ws.Range("D14").Value = "1.923.17"
$ws.Range("E14").Value = "  +1.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.85"
$ws.Range("E15").Value = "  +2.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.211"
$ws.Range("E16").Value = "  +5.36%  "

$ws.Range("D17").Value = "29.812.50"
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.96"
$ws.Range("E18").Value = "  +1.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.94"
$ws.Range("E19").Value = "  +2.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007816"
$ws.Range("E20").Value = "  +2.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.117"
$ws.Range("E22").Value = "  -0.83%  "

$ws.Range("D23").Value = "2.097.46"
$ws.Range("E23").Value = "  -2.25%  "

$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1587"
$ws.Range("E25").Value = "  -2.81%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.431"
$ws.Range("E26").Value = "  +1.80%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.75"
$ws.Range("E27").Value = "  -0.78%  "

$ws.Range("E28").Value = "  +0.71%  "

$ws.Range("E29").Value = "  -0.71%  "

$ws.Range("E30").Value = "  +2.71%  "

$ws.Range("E31").Value = "  +0.73%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.481"
$ws.Range("E32").Value = "  +2.51%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05568"
$ws.Range("E33").Value = "  -0.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.054"
$ws.Range("E34").Value = "  +0.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.238"
$ws.Range("E35").Value = "  -1.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7515"
$ws.Range("E36").Value = "  +2.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9972"
$ws.Range("E37").Value = "  -1.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.683"
$ws.Range("E38").Value = "  +1.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01937"
$ws.Range("E39").Value = "  +2.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.795"
$ws.Range("E40").Value = "  +0.86%  "

$ws.Range("D41").Value = "1.140.92"
$ws.Range("E41").Value = "  +12.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.33"
$ws.Range("E42").Value = "  +3.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4457"
$ws.Range("E43").Value = "  +1.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.954"
$ws.Range("E44").Value = "  +2.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8519"
$ws.Range("E45").Value = "  +2.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9998"
$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.888"
$ws.Range("E47").Value = "  +2.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.30"
$ws.Range("E48").Value = "  +0.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.054"
$ws.Range("E49").Value = "  +5.83%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.768"
$ws.Range("E50").Value = "  -0.53%  "

$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.511"
$ws.Range("E51").Value = "  +2.05%  "
